# ---------------------------------------------------------------------------
# Applies two changes to the deck:
#
#  1. The table on slide 16 gets its table style switched from the custom
#     "Table_0" style to the built-in style {B73A2977-BACA-4356-A7C9-88726AFA4EC1}.
#
#  2. The colour values that make up the presentation's theme colour scheme
#     (normally persisted as ppt/theme/theme1.xml, the theme used by the
#     slide master / all slides) are switched from the "Integral" palette to
#     the standard Office palette, via ThemeColorScheme so the written
#     <a:clrScheme> entries match what the document should contain.
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1. Table style -> built in style ------------------------------------
$slide = $p.Slides.Item(16)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{B73A2977-BACA-4356-A7C9-88726AFA4EC1}")
    }
}

# --- 2. Theme colour scheme -> Office palette ------------------------------
function BGR($r, $g, $b) {
    return ($b * 65536) + ($g * 256) + $r
}

$master = $p.Slides.Item(1).Master
$colorScheme = $master.Theme.ThemeColorScheme

# Index order matches the OOXML <a:clrScheme> child order:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink
$officePalette = @(
    (BGR 0x00 0x00 0x00),  # dk1
    (BGR 0xFF 0xFF 0xFF),  # lt1
    (BGR 0x44 0x54 0x6A),  # dk2
    (BGR 0xE7 0xE6 0xE6),  # lt2
    (BGR 0x5B 0x9B 0xD5),  # accent1
    (BGR 0xED 0x7D 0x31),  # accent2
    (BGR 0xA5 0xA5 0xA5),  # accent3
    (BGR 0xFF 0xC0 0x00),  # accent4
    (BGR 0x44 0x72 0xC4),  # accent5
    (BGR 0x70 0xAD 0x47),  # accent6
    (BGR 0x05 0x63 0xC1),  # hlink
    (BGR 0x95 0x4F 0x72)   # folHlink
)

for ($i = 1; $i -le 12; $i++) {
    $colorScheme.Item($i).RGB = $officePalette[$i - 1]
}
